$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "blend" label in D9 (no longer used)
$ws.Range("D9").ClearContents()

# 2. Row 19: "gauchissement" -> "distortion", "range or off" -> "range"
$ws.Range("B19").Value = "distortion"
$ws.Range("C19").Value = "range"

# 3. B22 ("wawa") gets vertical-centered alignment
$ws.Range("B22").VerticalAlignment = -4108

# 4. New row 23: continuation of the "wawa" block (merged with B22), with an extra "range" row
$ws.Range("A23").Value = 22
$ws.Range("C23").Value = "range"
$ws.Range("B22:B23").Merge()

# 5. New row 24: another "distortion"/"on/off" entry
$ws.Range("A24").Value = 23
$ws.Range("A24").HorizontalAlignment = -4108
$ws.Range("A24").VerticalAlignment = -4108
$ws.Range("B24").Value = "distortion"
$ws.Range("C24").Value = "on/off"

# 6. Update selection to match the authored state
[void]$ws.Range("E19").Select()
